# "Dados Obrigatórios e Não Obrigatórios"
#
# Módulos!B6 ("calcular_acoes_regressivas_inss") flips from Obrigatório (TRUE)
# to Não Obrigatório (FALSE), and Módulos!B7 ("calcular_beneficios_inss")
# flips the other way, from FALSE to TRUE. Funcoes_Inputs!D31:D61 are
# VLOOKUP(...)-driven off that table, so they recalc automatically once the
# source flags change. Also widen column A on Módulos (the longer label no
# longer "best fits" at the old width) and leave the selection cursor sitting
# on Módulos!B6 - the cell that was just edited - while Funcoes_Inputs stays
# the active tab.

$wb = $excel.ActiveWorkbook

$funcoes = $wb.Worksheets.Item("Funcoes_Inputs")
$modulos = $wb.Worksheets.Item("Módulos")

$modulos.Activate()

# Column A was auto-fit (bestFit) before; now it's a manually-set width.
$modulos.Columns.Item(1).ColumnWidth = 32

# Swap which module is flagged mandatory.
$modulos.Range("B6").Formula = "=FALSE"
$modulos.Range("B7").Formula = "=TRUE"

$modulos.Range("B6").Select()

# Return focus to the sheet that stays active in the saved workbook.
$funcoes.Activate()
$funcoes.Range("A1").Select()
